$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'272.09"
$ws.Range("G2").Value = "'18"
$ws.Range("D3").Value = "'23.08"
$ws.Range("G3").Value = "'18"
$ws.Range("D4").Value = "'6.355"
$ws.Range("G4").Value = "'18"
$ws.Range("D5").Value = "'0.06297"
$ws.Range("G5").Value = "'18"
$ws.Range("D6").Value = "'3.664"
$ws.Range("G6").Value = "'18"
$ws.Range("D7").Value = "'6.756"
$ws.Range("G7").Value = "'18"
$ws.Range("D8").Value = "'1.398"
$ws.Range("G8").Value = "'18"
$ws.Range("D9").Value = "'0.8370"
$ws.Range("G9").Value = "'18"
$ws.Range("D10").Value = "'0.1627"
$ws.Range("G10").Value = "'18"
$ws.Range("D11").Value = "'0.08391"
$ws.Range("G11").Value = "'18"
$ws.Range("D12").Value = "'0.03445"
$ws.Range("G12").Value = "'18"
$ws.Range("D13").Value = "'0.03194"
$ws.Range("G13").Value = "'18"
$ws.Range("D14").Value = "'0.09309"
$ws.Range("G14").Value = "'18"
$ws.Range("D15").Value = "'3.927"
$ws.Range("G15").Value = "'18"
$ws.Range("D16").Value = "'0.001717"
$ws.Range("G16").Value = "'18"
$ws.Range("D17").Value = "'0.04861"
$ws.Range("G17").Value = "'18"
$ws.Range("D18").Value = "'0.006283"
$ws.Range("G18").Value = "'18"
$ws.Range("D19").Value = "'0.005479"
$ws.Range("G19").Value = "'18"
$ws.Range("G20").Value = "'18"
$ws.Range("D21").Value = "'0.0001499"
$ws.Range("G21").Value = "'18"
$ws.Range("D22").Value = "'3.728"
$ws.Range("G22").Value = "'18"
$ws.Range("D23").Value = "'2.347"
$ws.Range("G23").Value = "'18"
$ws.Range("D24").Value = "'0.01383"
$ws.Range("G24").Value = "'18"
$ws.Range("D25").Value = "'0.3377"
$ws.Range("G25").Value = "'18"
$ws.Range("D26").Value = "'0.1218"
$ws.Range("G26").Value = "'18"
$ws.Range("D27").Value = "'0.0002680"
$ws.Range("G27").Value = "'18"
$ws.Range("G28").Value = "'18"
$ws.Range("G29").Value = "'18"
$ws.Range("G30").Value = "'18"
$ws.Range("G31").Value = "'18"
$ws.Range("G32").Value = "'18"
$ws.Range("G33").Value = "'18"
$ws.Range("G34").Value = "'18"
$ws.Range("G35").Value = "'18"
$ws.Range("G36").Value = "'18"
$ws.Range("G37").Value = "'18"
$ws.Range("G38").Value = "'18"
$ws.Range("G39").Value = "'18"
$ws.Range("D40").Value = "'0.04696"
$ws.Range("G40").Value = "'18"
$ws.Range("D41").Value = "'0.006891"
$ws.Range("G41").Value = "'18"
$ws.Range("D42").Value = "'0.1175"
$ws.Range("G42").Value = "'18"
$ws.Range("D43").Value = "'0.003458"
$ws.Range("G43").Value = "'18"
$ws.Range("D44").Value = "'0.01260"
$ws.Range("G44").Value = "'18"
$ws.Range("D45").Value = "'0.00006243"
$ws.Range("G45").Value = "'18"
$ws.Range("G46").Value = "'18"
$ws.Range("D47").Value = "'0.6994"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("G47").Value = "'18"
$ws.Range("D48").Value = "'0.1199"
$ws.Range("G48").Value = "'18"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "48CryptobidCoinCBC"
$ws.Range("G49").Value = "'18"
$ws.Range("D50").Value = "'0.01239"
$ws.Range("G50").Value = "'18"
$ws.Range("G51").Value = "'18"
